$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "23.753.52"
$ws.Range("E2").Value = "  -3.11%  "
$ws.Range("D3").Value = "1.611.24"
$ws.Range("E3").Value = "  -3.47%  "
Set-TextValue $ws.Range("D4") "1.007"
$ws.Range("E4").Value = "  +0.35%  "
Set-TextValue $ws.Range("D5") "1.004"
$ws.Range("E5").Value = "  +0.11%  "
Set-TextValue $ws.Range("D6") "305.59"
$ws.Range("E6").Value = "  -2.62%  "
Set-TextValue $ws.Range("D7") "0.3889"
$ws.Range("E7").Value = "  -0.31%  "
Set-TextValue $ws.Range("D8") "0.3807"
$ws.Range("E8").Value = "  -3.06%  "
Set-TextValue $ws.Range("D9") "1.004"
$ws.Range("E9").Value = "  +0.10%  "
Set-TextValue $ws.Range("D10") "1.343"
$ws.Range("E10").Value = "  -4.62%  "
Set-TextValue $ws.Range("D11") "48.62"
$ws.Range("E11").Value = "  -5.50%  "
Set-TextValue $ws.Range("D12") "0.08380"
$ws.Range("E12").Value = "  -2.72%  "
Set-TextValue $ws.Range("D13") "23.54"
$ws.Range("E13").Value = "  -5.78%  "
Set-TextValue $ws.Range("D14") "6.927"
$ws.Range("E14").Value = "  -4.77%  "
Set-TextValue $ws.Range("D15") "0.00001265"
$ws.Range("E15").Value = "  -3.49%  "
Set-TextValue $ws.Range("D16") "7.376"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("D17").Value = "1.610.48"
$ws.Range("E17").Value = "  -7.30%  "
Set-TextValue $ws.Range("D18") "92.87"
$ws.Range("E18").Value = "  -0.48%  "
Set-TextValue $ws.Range("D19") "0.06915"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("E20").Value = "  -3.58%  "
Set-TextValue $ws.Range("D21") "6.764"
$ws.Range("E21").Value = "  -3.94%  "
Set-TextValue $ws.Range("D22") "1.004"
$ws.Range("E22").Value = "  -0.14%  "
Set-TextValue $ws.Range("D23") "13.30"
$ws.Range("E23").Value = "  -4.65%  "
$ws.Range("D24").Value = "23.796.08"
$ws.Range("E24").Value = "  -2.92%  "
Set-TextValue $ws.Range("D25") "2.413"
$ws.Range("E25").Value = "  +1.43%  "
Set-TextValue $ws.Range("D26") "2.724"
$ws.Range("E26").Value = "  -0.20%  "
Set-TextValue $ws.Range("D27") "21.97"
$ws.Range("E27").Value = "  -5.02%  "
$ws.Range("E28").Value = "  -2.45%  "
Set-TextValue $ws.Range("D29") "139.19"
$ws.Range("E29").Value = "  -5.41%  "
Set-TextValue $ws.Range("D30") "5.263"
$ws.Range("E30").Value = "  -8.94%  "
Set-TextValue $ws.Range("D31") "7.723"
$ws.Range("E31").Value = "  -6.86%  "
Set-TextValue $ws.Range("D32") "2.459"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").Value = "1.787.88"
$ws.Range("E33").Value = "  -6.75%  "
Set-TextValue $ws.Range("D34") "0.07998"
$ws.Range("E34").Value = "  -3.74%  "
Set-TextValue $ws.Range("D35") "0.9572"
$ws.Range("E35").Value = "  -1.86%  "
Set-TextValue $ws.Range("D38") "0.2639"
$ws.Range("E38").Value = "  -5.45%  "
Set-TextValue $ws.Range("D39") "0.09106"
$ws.Range("E39").Value = "  -3.57%  "
Set-TextValue $ws.Range("D40") "10.24"
$ws.Range("E40").Value = "  -0.36%  "
Set-TextValue $ws.Range("D41") "13.21"
$ws.Range("E41").Value = "  -1.91%  "
Set-TextValue $ws.Range("D42") "1.412"
$ws.Range("E42").Value = "  -8.36%  "
Set-TextValue $ws.Range("D43") "0.7389"
$ws.Range("E43").Value = "  -5.93%  "
Set-TextValue $ws.Range("D44") "15.87"
$ws.Range("E44").Value = "  -3.28%  "
Set-TextValue $ws.Range("D45") "0.6790"
$ws.Range("E45").Value = "  -4.05%  "
Set-TextValue $ws.Range("D46") "2.415"
$ws.Range("E46").Value = "  -5.05%  "
Set-TextValue $ws.Range("D47") "4.042"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("E48").Value = "  +0.06%  "
Set-TextValue $ws.Range("D49") "0.08197"
$ws.Range("E49").Value = "  -4.33%  "
Set-TextValue $ws.Range("D50") "131.61"
$ws.Range("E50").Value = "  -3.96%  "
Set-TextValue $ws.Range("D51") "1.235"
$ws.Range("E51").Value = "  -6.46%  "

# Row 36/37 swap: VeChain <-> InternetComputer(DFINITY), with updated data
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D36") "6.568"
$ws.Range("E36").Value = "  -5.30%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02842"
$ws.Range("E37").Value = "  -5.92%  "
